$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.335.90'
$ws.Range('E2').Value = '  -1.86%  '

# Row 3
$ws.Range('D3').Value = '2.635.42'
$ws.Range('E3').Value = '  +0.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').Value = '514.97'
$ws.Range('E5').Value = '  -1.42%  '

# Row 6
$ws.Range('D6').Value = '148.53'
$ws.Range('E6').Value = '  -2.57%  '

# Row 7
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.25%  '

# Row 8
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('D9').Value = '2.666.81'
$ws.Range('E9').Value = '  +1.28%  '

# Row 10
$ws.Range('E10').Value = '  +0.87%  '

# Row 11
$ws.Range('D11').Value = '0.106'
$ws.Range('E11').Value = '  -1.34%  '

# Row 12
$ws.Range('D12').Value = '0.341'
$ws.Range('E12').Value = '  -1.45%  '

# Row 13
$ws.Range('E13').Value = '  -0.91%  '

# Row 14
$ws.Range('D14').Value = '3.100.61'
$ws.Range('E14').Value = '  +0.36%  '

# Row 15
$ws.Range('D15').Value = '59.231.01'
$ws.Range('E15').Value = '  -2.03%  '

# Row 16
$ws.Range('D16').Value = '21.44'
$ws.Range('E16').Value = '  -0.81%  '

# Row 17
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -0.99%  '

# Row 18
$ws.Range('D18').Value = '2.663.21'
$ws.Range('E18').Value = '  +1.20%  '

# Row 19
$ws.Range('D19').Value = '4.62'
$ws.Range('E19').Value = '  -1.32%  '

# Row 20
$ws.Range('D20').Value = '345.67'
$ws.Range('E20').Value = '  -0.79%  '

# Row 21
$ws.Range('D21').Value = '10.54'
$ws.Range('E21').Value = '  +0.34%  '

# Row 22
$ws.Range('D22').Value = '6.21'
$ws.Range('E22').Value = '  -0.15%  '

# Row 23
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  +0.29%  '

# Row 24
$ws.Range('D24').Value = '61.16'
$ws.Range('E24').Value = '  +0.39%  '

# Row 25
$ws.Range('D25').Value = '0.427'
$ws.Range('E25').Value = '  +0.61%  '

# Row 26
$ws.Range('D26').Value = '2.759.65'
$ws.Range('E26').Value = '  +0.30%  '

# Row 27
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.994'
$ws.Range('E27').Value = '  -0.63%  '

# Row 28
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = '0.161'
$ws.Range('E28').Value = '  -2.47%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0836'
$ws.Range('E29').Value = '  -0.52%  '

# Row 30
$ws.Range('D30').Value = '7.16'
$ws.Range('E30').Value = '  -0.22%  '

# Row 31
$ws.Range('D31').Value = '6.61'
$ws.Range('E31').Value = '  +8.22%  '

# Row 32
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.22%  '

# Row 33
$ws.Range('D33').Value = '19.04'
$ws.Range('E33').Value = '  -0.44%  '

# Row 34
$ws.Range('E34').Value = '  -1.76%  '

# Row 35
$ws.Range('D35').Value = '149.56'
$ws.Range('E35').Value = '  -0.11%  '

# Row 36
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +19.34%  '

# Row 37
$ws.Range('D37').Value = '4.07'
$ws.Range('E37').Value = '  +0.55%  '

# Row 38
$ws.Range('D38').Value = '1.16'
$ws.Range('E38').Value = '  -1.02%  '

# Row 39
$ws.Range('D39').Value = '0.878'
$ws.Range('E39').Value = '  -2.06%  '

# Row 40
$ws.Range('D40').Value = '36.51'
$ws.Range('E40').Value = '  -0.07%  '

# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '3.72'
$ws.Range('E41').Value = '  +0.74%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.44'
$ws.Range('E42').Value = '  -0.96%  '

# Row 43
$ws.Range('D43').Value = '288.02'
$ws.Range('E43').Value = '  -5.34%  '

# Row 44
$ws.Range('D44').Value = '0.629'
$ws.Range('E44').Value = '  -0.36%  '

# Row 45
$ws.Range('D45').Value = '0.0999'
$ws.Range('E45').Value = '  -0.82%  '

# Row 46
$ws.Range('D46').Value = '0.993'
$ws.Range('E46').Value = '  -0.46%  '

# Row 47
$ws.Range('D47').Value = '19.73'
$ws.Range('E47').Value = '  -0.20%  '

# Row 48
$ws.Range('D48').Value = '0.0546'
$ws.Range('E48').Value = '  -1.77%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.80'
$ws.Range('E49').Value = '  -0.31%  '

# Row 50
$ws.Range('D50').Value = '0.0233'
$ws.Range('E50').Value = '  -1.81%  '

# Row 51
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '10.26'
$ws.Range('E51').Value = '  -0.97%  '
